$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3137.375
$ws.Range("I8").Value = 20.2
$ws.Range("J8").Value = 8332.666999999999
$ws.Range("K8").Value = 60.59999999999999
$ws.Range("L8").Value = 24998.001
$ws.Range("M8").Value = 78.40000000000001
$ws.Range("N8").Value = -25276.001

$ws.Range("H11").Value = 2117
$ws.Range("I11").Value = 2117
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2117
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1977

$ws.Range("H32").Value = 999.5
$ws.Range("I32").Value = 999.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 999.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -673.5

$ws.Range("H33").Value = 155.125
$ws.Range("I33").Value = 155.125
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 155.125
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 73.875

$ws.Range("H42").Value = 53.125
$ws.Range("I42").Value = 55
$ws.Range("J42").Value = 50
$ws.Range("K42").Value = 165
$ws.Range("L42").Value = 150
$ws.Range("M42").Value = 65

$ws.Range("H47").Value = 4500
$ws.Range("I47").Value = 4500
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 4500
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -3528

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082
$ws.Range("M100").ClearContents()

$ws.Range("H115").Value = 1425
$ws.Range("I115").Value = 850
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 2550
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -983
$ws.Range("N115").Value = -9134

$ws.Range("H123").Value = 200780
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 200780
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 200780
$ws.Range("N123").Value = -210580

$ws.Range("H136").Value = 89780
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 89780
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 89780
$ws.Range("N136").Value = -99980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7583.4
$ws.Range("I61").Value = 4690.5713
$ws.Range("J61").Value = 14333.333
$ws.Range("K61").Value = 4690.5713
$ws.Range("L61").Value = 14333.333
$ws.Range("M61").Value = -4478.5713

$ws.Range("H63").Value = 1995
$ws.Range("I63").Value = 1995
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1995
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1309

$ws.Range("H66").Value = 1995
$ws.Range("I66").Value = 1995
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9975
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6543

$ws.Range("H97").Value = 2010.5
$ws.Range("I97").Value = 2010
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 2010
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -1514
$ws.Range("N97").Value = -3003

$ws.Range("H102").Value = 1345.1428
$ws.Range("I102").Value = 1345.1428
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1345.1428
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 276.8571999999999

$ws.Range("H132").Value = 9969.362999999999
$ws.Range("I132").Value = 3667.5715
$ws.Range("J132").Value = 20997.5
$ws.Range("K132").Value = 11002.7145
$ws.Range("L132").Value = 62992.5
$ws.Range("M132").Value = -8472.7145

$ws.Range("H136").Value = 7583.4
$ws.Range("I136").Value = 4690.5713
$ws.Range("J136").Value = 14333.333
$ws.Range("K136").Value = 14071.7139
$ws.Range("L136").Value = 42999.999
$ws.Range("M136").Value = -11521.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 201.33333
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 227
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 227
$ws.Range("M7").Value = -37

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -253

$ws.Range("H134").Value = 8263
$ws.Range("I134").Value = 2568.2
$ws.Range("J134").Value = 22500
$ws.Range("K134").Value = 7704.599999999999
$ws.Range("L134").Value = 67500
$ws.Range("M134").Value = -5169.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 14618.5
$ws.Range("I19").Value = 237.5
$ws.Range("J19").Value = 28999.5
$ws.Range("K19").Value = 237.5
$ws.Range("L19").Value = 28999.5
$ws.Range("M19").Value = -67.5
$ws.Range("N19").Value = -29339.5

$ws.Range("H24").Value = 14618.5
$ws.Range("I24").Value = 237.5
$ws.Range("J24").Value = 28999.5
$ws.Range("K24").Value = 237.5
$ws.Range("L24").Value = 28999.5
$ws.Range("M24").Value = -67.5
$ws.Range("N24").Value = -29339.5

$ws.Range("H41").Value = 19999
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19999
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19999
$ws.Range("N41").Value = -20855

$ws.Range("H58").Value = 13499.333
$ws.Range("I58").Value = 3999.5
$ws.Range("J58").Value = 18249.25
$ws.Range("K58").Value = 3999.5
$ws.Range("L58").Value = 18249.25
$ws.Range("M58").Value = -3796.5
$ws.Range("N58").Value = -18655.25

$ws.Range("H100").Value = 50000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 50000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164

$ws.Range("H115").Value = 49290
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 49290
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 49290
$ws.Range("N115").Value = -51640

$ws.Range("H132").Value = 6030.846
$ws.Range("I132").Value = 4522.3335
$ws.Range("J132").Value = 9425
$ws.Range("K132").Value = 13567.0005
$ws.Range("L132").Value = 28275
$ws.Range("M132").Value = -11037.0005

$ws.Range("H134").Value = 8255.571
$ws.Range("I134").Value = 2322.25
$ws.Range("J134").Value = 16166.667
$ws.Range("K134").Value = 6966.75
$ws.Range("L134").Value = 48500.001
$ws.Range("M134").Value = -4431.75

$ws.Range("H136").Value = 13499.333
$ws.Range("I136").Value = 3999.5
$ws.Range("J136").Value = 18249.25
$ws.Range("K136").Value = 11998.5
$ws.Range("L136").Value = 54747.75
$ws.Range("M136").Value = -9448.5
$ws.Range("N136").Value = -59847.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 76930340
$ws.Range("I4").Value = 350
$ws.Range("J4").Value = 111121450
$ws.Range("K4").Value = 1050
$ws.Range("L4").Value = 333364350
$ws.Range("M4").Value = -938
$ws.Range("N4").Value = -333364574

$ws.Range("H7").Value = 27.333334
$ws.Range("I7").Value = 30
$ws.Range("J7").Value = 14
$ws.Range("K7").Value = 90
$ws.Range("L7").Value = 42
$ws.Range("M7").Value = 22

$ws.Range("H44").Value = 539.75
$ws.Range("I44").Value = 553
$ws.Range("J44").Value = 500
$ws.Range("K44").Value = 1659
$ws.Range("L44").Value = 1500
$ws.Range("M44").Value = -1261
$ws.Range("N44").Value = -2296

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H76").Value = 100
$ws.Range("I76").Value = 100
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 300
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 83

$ws.Range("H79").Value = 100
$ws.Range("I79").Value = 100
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 300
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 1026

$ws.Range("H97").Value = 2076.5
$ws.Range("I97").Value = 1950
$ws.Range("J97").Value = 2203
$ws.Range("K97").Value = 5850
$ws.Range("L97").Value = 6609
$ws.Range("M97").Value = -5354
$ws.Range("N97").Value = -7601

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2730

$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2064

$ws.Range("H102").Value = 1001.6667
$ws.Range("I102").Value = 1001.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1001.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 620.3333

$ws.Range("H113").Value = 4583
$ws.Range("I113").Value = 4699.6
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 4699.6
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -2529.6
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1750
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1364

$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3298
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1330
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
